$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "VEhop GAMO, 24Hrs Play, 60ms Low Latency, On Ear Wireless With Mic Headphones/Earphones Black"
$ws.Range("B1").Value = "Rs. 1,199"

$ws.Range("A2").Value = "pTron Bassbuds Sports On Ear Bluetooth Headphone 48 Hours Playback IPX4(Splash & Sweat Proof) Passive noise cancellation -Bluetooth V 5.2 Black"
$ws.Range("B2").Value = "Rs. 999"

$ws.Range("A3").Value = "boAt Airdopes 131/138 Twin Wireless Earbuds with IWP Technology, Bluetooth V5.0, Immersive Audio, Up to 15H Total Playback, Instant Voice Assistant and Type-C Charging,Bluetooth Earphone (Active Black)"
$ws.Range("B3").Value = "Rs. 1,499"

$ws.Range("A4").Value = "Tecsox PowerHouse Earbud In Ear Bluetooth Earphone 45 Hours Playback Bluetooth IPX5(Splash Proof) Powerfull Bass -Bluetooth V 5.1 Black"
$ws.Range("B4").Value = "Rs. 725"

$ws.Range("A5").Value = "Lenovo QE03 In Ear Bluetooth Neckband 8 Hours Playback IPX5(Splash & Sweat Proof) Powerfull bass -Bluetooth V 5.0 Black"
$ws.Range("B5").Value = "Rs. 899"
